# Update countries & provincias Spain
# Refresh the "Pais" data pull: new case counts for a handful of countries,
# two countries swap rank (so their name labels trade rows), and the
# "datos actualizados" timestamp advances.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp banner (row 1)
$ws.Range("A1").Value = "Datos actualizados a 22 de Agosto de 2020 a las 08:51"

# Row 6 - India: refreshed totals
$ws.Range("B6").Value = 2979562
$ws.Range("C6").Value = 6194
$ws.Range("D6").Value = 2223202
$ws.Range("E6").Value = 700410
$ws.Range("G6").Value = 22
$ws.Range("H6").Value = 55950

# Rows 60/61 - Uzbekistan overtakes Venezuela in total cases, so the two
# countries swap places (row 60 becomes Uzbekistan with its new, higher
# totals; row 61 becomes Venezuela, keeping its previous totals).
$ws.Range("A60").Value = "Uzbekistan"
$ws.Range("B60").Value = 38231
$ws.Range("C60").Value = 157
$ws.Range("D60").Value = 33989
$ws.Range("E60").Value = 3980
$ws.Range("G60").Value = 2
$ws.Range("H60").Value = 262

$ws.Range("A61").Value = "Venezuela"
$ws.Range("B61").Value = 38219
$ws.Range("D61").Value = 27306
$ws.Range("E61").Value = 10596
$ws.Range("H61").Value = 317

# Row 72 - Australia: refreshed totals
$ws.Range("D72").Value = 18759
$ws.Range("E72").Value = 5358

# Row 73 - El Salvador: refreshed totals
$ws.Range("B73").Value = 24420
$ws.Range("C73").Value = 220
$ws.Range("D73").Value = 12021
$ws.Range("E73").Value = 11745
$ws.Range("G73").Value = 8
$ws.Range("H73").Value = 654

# Row 149 - Georgia: refreshed totals
$ws.Range("B149").Value = 1394
$ws.Range("C149").Value = 9
$ws.Range("D149").Value = 1132
$ws.Range("E149").Value = 245

# Rows 202/203 - Timor Oriental ties/overtakes Santa Lucia, so the two
# countries swap places (row 202 becomes Timor Oriental; row 203 becomes
# Santa Lucia, keeping its previous totals).
$ws.Range("A202").Value = "Timor Oriental"
$ws.Range("A203").Value = "Santa Lucia"
$ws.Range("B203").Value = 26
$ws.Range("E203").Value = 1
